# A new weekly price record (2023-08-04) was inserted into the logged
# series at sheet row 708, pushing every existing row from 708 downward
# down by one (last row moves from 798 to 799).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 708 - this shifts rows 708:798 down to
# 709:799 and carries the date-format style of the surrounding rows.
$ws.Rows("708:708").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A708").Value = 6
$ws.Range("B708").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C708").Value = "Metropolitana"
$ws.Range("D708").Value = 45142
$ws.Range("E708").Value = 13
$ws.Range("F708").Value = 100112012
$ws.Range("G708").Value = "Espinaca"
$ws.Range("H708").Value = "Sin especificar"
$ws.Range("I708").Value = "Primera"
$ws.Range("J708").Value = 600
$ws.Range("K708").Value = 5000
$ws.Range("L708").Value = 6000
$ws.Range("M708").Value = 5417
$ws.Range("N708").Value = "`$/cuna 10 kilos"
$ws.Range("O708").Value = "Región Metropolitana"
$ws.Range("P708").Value = 542
$ws.Range("Q708").Value = 10
$ws.Range("R708").Value = "Hortaliza"
